# New tenant support in live
# Appends new sprint-run rows to the AMSIN, BETA and AMS history sheets.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Add-HistoryRow($ws, $row, $dateText, $timeSerial, $sprintName, $total, $pass, $fail, $taken, $textFormatSourceRange, $dateFormatSourceRange) {
    # Column A holds a plain date-looking label, but it must stay literal
    # text (matching the source data) rather than be auto-parsed into a
    # real Excel date serial.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $dateText
    # Put the cell format back to a plain/general style (matches the other
    # text cells in the column) now that the literal text is committed.
    $textFormatSourceRange.Copy()
    $ws.Cells.Item($row, 1).PasteSpecial($xlPasteFormats)

    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $taken

    # Column B is a real date/time serial; copy the custom date-time
    # number format from an existing dated cell, then write the value.
    $dateFormatSourceRange.Copy()
    $ws.Cells.Item($row, 2).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($row, 2).Value = $timeSerial
}

# ---------------------------------------------------------------------------
# AMSIN sheet: add rows 65 and 66. Row 64 (the previous last row) also picks
# up a tiny refreshed run-time value.
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Range("B64").Value = 45034.55313322916

Add-HistoryRow $wsAmsin 65 "2023-04-19" 45035.70167424768 "176scndcyc" 105 94 11 5.35 $wsAmsin.Range("A63") $wsAmsin.Range("B63")
Add-HistoryRow $wsAmsin 66 "2023-04-20" 45036.42445592592 "176fnlruntest" 105 105 0 3.69 $wsAmsin.Range("A63") $wsAmsin.Range("B63")

# ---------------------------------------------------------------------------
# BETA sheet: add row 32.
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Add-HistoryRow $wsBeta 32 "2023-04-20" 45036.52280033565 "176beta" 105 105 0 2.72 $wsBeta.Range("A31") $wsBeta.Range("B31")

# ---------------------------------------------------------------------------
# AMS sheet: add rows 33 and 34.
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

Add-HistoryRow $wsAms 33 "2023-05-04" 45050.72727009259 "176firstsycle" 105 105 0 2.92 $wsAms.Range("A32") $wsAms.Range("B32")
Add-HistoryRow $wsAms 34 "2023-05-08" 45054.5429728984 "176htfxtrl" 105 105 0 3.05 $wsAms.Range("A32") $wsAms.Range("B32")
